$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# Update the plant_type data_value cell from "plants" to "planting"
$survey.Range("D5").Value = "planting"

# Update selections on each sheet
[void]$survey.Range("D5").Select()
[void]$settings.Range("B15").Select()

# Make the survey sheet the active tab
[void]$survey.Activate()
